$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 10786.9707726745
    "C2" = 9955.7575759944
    "E2" = 5554.28084147527
    "F2" = 177.109934061236

    "B3" = 10562.1689847619
    "C3" = 9458.68939678446
    "E3" = 6155.73199951375
    "F3" = 275.459224845759

    "B4" = 4534.81818341195
    "C4" = 7304.09929676796
    "E4" = 6133.34423387248
    "F4" = 184.751813776685

    "B5" = 4626.3092054951
    "C5" = 7649.51030094833
    "E5" = 6416.61534177313
    "F5" = 210.946901780061

    "B6" = 12139.8470618849
    "C6" = 10883.9297484406
    "E6" = 7218.64330302934
    "F6" = 379.132210477915

    "B7" = 11902.6345140631
    "C7" = 10701.1707348469
    "E7" = 7002.123017967
    "F7" = 362.495573033914

    "C9" = 10379.4632813111
    "F9" = 324.387652447176

    "C10" = 10060.9430782958
    "F10" = 311.115977321541

    "C11" = 7178.91858826329
    "F11" = 175.355774689018

    "C12" = 6949.13466041134
    "F12" = 165.427753446778

    "C13" = 10542.5891170901
    "E13" = 7489.86308189519
    "F13" = 390.341341624388

    "C14" = 10311.0206618965
    "E14" = 7489.86308189519
    "F14" = 380.692655991322

    "C15" = 10228.9567679552
    "E15" = 7489.86308189519
    "F15" = 377.2733270771
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
